# Update the GPA figure (a typo fix for the Spring 2017 semester) and
# move the auto "_GoBack" bookmark so it again sits right after the last
# edit location (immediately after the GPA value), matching what Word
# does automatically whenever a document is resaved after an edit.

$d = $word.ActiveDocument

# 1. Fix the GPA value: 3.656 -> 3.636
$d.Content.Find.Execute("3.656", $true, $false, $false, $false, $false,
                         $true, 1, $false, "3.636", 2)

# 2. Locate the end of the run we just edited (right after "3.636",
#    still inside its paragraph, before the paragraph mark).
$found = $d.Content
$found.Find.Execute("3.636", $true, $false, $false, $false, $false,
                     $true, 1, $false, "", 0)
$insertPos = $found.End

# A zero-length Range sitting exactly on a paragraph-end boundary isn't
# addressed reliably, so nudge past it: insert a throwaway character
# right after the GPA text, drop the bookmark there, then remove the
# throwaway character again. This leaves the bookmark collapsed right
# after "3.636" and before the paragraph mark, exactly where it needs
# to be.
$anchor = $d.Range($insertPos, $insertPos)
$anchor.InsertAfter("X")

$bmRange = $d.Range($insertPos, $insertPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$scratch = $d.Range($insertPos, $insertPos + 1)
$scratch.Delete()
